$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.072.94'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.931.50'
$ws.Range("E3").Value = '  +3.24%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'604.06"
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").Value = "'167.89"
$ws.Range("E6").Value = '  +2.38%  '
$ws.Range("D7").Value = '3.930.57'
$ws.Range("E7").Value = '  +3.28%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = "'0.173"
$ws.Range("E10").Value = '  +2.76%  '
$ws.Range("D11").Value = "'6.47"
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("D12").Value = "'0.466"
$ws.Range("E12").Value = '  +1.60%  '
$ws.Range("D13").Value = "'0.0000260"
$ws.Range("E13").Value = '  +6.42%  '
$ws.Range("D14").Value = "'37.63"
$ws.Range("E14").Value = '  +1.82%  '
$ws.Range("D15").Value = '4.596.79'
$ws.Range("E15").Value = '  +3.40%  '
$ws.Range("D16").Value = '3.907.61'
$ws.Range("E16").Value = '  +1.75%  '
$ws.Range("D17").Value = '69.144.48'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = "'7.49"
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = "'17.40"
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("E20").Value = '  -1.84%  '
$ws.Range("D21").Value = "'11.00"
$ws.Range("E21").Value = '  -4.62%  '
$ws.Range("D22").Value = "'496.43"
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("D23").Value = "'0.732"
$ws.Range("E23").Value = '  +1.89%  '
$ws.Range("D24").Value = "'0.0000169"
$ws.Range("E24").Value = '  +7.93%  '
$ws.Range("D25").Value = "'85.16"
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("E26").Value = '  +1.72%  '
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("E28").Value = '  +2.31%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("D31").Value = '4.091.05'
$ws.Range("E31").Value = '  +3.23%  '
$ws.Range("D32").Value = "'2.39"
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("D33").Value = "'7.79"
$ws.Range("E33").Value = '  -2.33%  '
$ws.Range("D34").Value = "'32.00"
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("D35").Value = '3.904.83'
$ws.Range("E35").Value = '  +4.09%  '
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").Value = "'1.04"
$ws.Range("E37").Value = '  +1.64%  '
$ws.Range("D38").Value = "'6.01"
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = "'3.33"
$ws.Range("E39").Value = '  +10.00%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = "'0.139"
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  +1.31%  '
$ws.Range("D43").Value = "'433.79"
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("E44").Value = '  +1.65%  '
$ws.Range("D45").Value = "'47.95"
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").Value = "'8.61"
$ws.Range("E46").Value = '  +3.04%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = "'0.000287"
$ws.Range("E48").Value = '  +27.54%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = "'0.0366"
$ws.Range("E49").Value = '  +3.46%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = "'143.08"
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.808.76'
$ws.Range("E51").Value = '  -0.47%  '
